# Apply crypto price/volume updates for Sat Aug 12 08:12:39 UTC 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell="D2"; Value='29.400.72'; ForceText=$false},
    @{Cell="D3"; Value='1.849.19'; ForceText=$false},
    @{Cell="E3"; Value='  +0.07%  '; ForceText=$false},
    @{Cell="D4"; Value='0.9995'; ForceText=$true},
    @{Cell="E4"; Value='  +0.09%  '; ForceText=$false},
    @{Cell="D5"; Value='240.32'; ForceText=$true},
    @{Cell="E5"; Value='  -0.07%  '; ForceText=$false},
    @{Cell="D6"; Value='0.6301'; ForceText=$true},
    @{Cell="E6"; Value='  -0.05%  '; ForceText=$false},
    @{Cell="D7"; Value='1.000'; ForceText=$true},
    @{Cell="E7"; Value='  +0.03%  '; ForceText=$false},
    @{Cell="D8"; Value='0.07621'; ForceText=$true},
    @{Cell="E8"; Value='  +0.93%  '; ForceText=$false},
    @{Cell="D9"; Value='0.2936'; ForceText=$true},
    @{Cell="E9"; Value='  -0.77%  '; ForceText=$false},
    @{Cell="D10"; Value='24.53'; ForceText=$true},
    @{Cell="E10"; Value='  +0.23%  '; ForceText=$false},
    @{Cell="D11"; Value='0.07745'; ForceText=$true},
    @{Cell="E11"; Value='  +0.38%  '; ForceText=$false},
    @{Cell="D12"; Value='1.840.82'; ForceText=$false},
    @{Cell="E12"; Value='  -0.32%  '; ForceText=$false},
    @{Cell="D13"; Value='5.008'; ForceText=$true},
    @{Cell="E13"; Value='  +0.22%  '; ForceText=$false},
    @{Cell="D14"; Value='0.00001087'; ForceText=$true},
    @{Cell="E14"; Value='  +8.37%  '; ForceText=$false},
    @{Cell="D15"; Value='0.6794'; ForceText=$true},
    @{Cell="E15"; Value='  -0.85%  '; ForceText=$false},
    @{Cell="D16"; Value='83.54'; ForceText=$true},
    @{Cell="E16"; Value='  +0.58%  '; ForceText=$false},
    @{Cell="D17"; Value='2.086.00'; ForceText=$false},
    @{Cell="E17"; Value='  -7.88%  '; ForceText=$false},
    @{Cell="D18"; Value='6.149'; ForceText=$true},
    @{Cell="E18"; Value='  -0.05%  '; ForceText=$false},
    @{Cell="D19"; Value='29.428.56'; ForceText=$false},
    @{Cell="D20"; Value='228.85'; ForceText=$true},
    @{Cell="E20"; Value='  -0.41%  '; ForceText=$false},
    @{Cell="E21"; Value='  +0.08%  '; ForceText=$false},
    @{Cell="E22"; Value='  +0.03%  '; ForceText=$false},
    @{Cell="D23"; Value='7.459'; ForceText=$true},
    @{Cell="E23"; Value='  -1.41%  '; ForceText=$false},
    @{Cell="D24"; Value='1.001'; ForceText=$true},
    @{Cell="E24"; Value='  +0.01%  '; ForceText=$false},
    @{Cell="D25"; Value='157.27'; ForceText=$true},
    @{Cell="E25"; Value='  +0.14%  '; ForceText=$false},
    @{Cell="D26"; Value='0.1389'; ForceText=$true},
    @{Cell="E26"; Value='  -0.79%  '; ForceText=$false},
    @{Cell="D27"; Value='8.371'; ForceText=$true},
    @{Cell="E27"; Value='  -0.09%  '; ForceText=$false},
    @{Cell="E28"; Value='  -0.09%  '; ForceText=$false},
    @{Cell="E29"; Value='  +0.08%  '; ForceText=$false},
    @{Cell="D30"; Value='1.309'; ForceText=$true},
    @{Cell="E30"; Value='  +4.35%  '; ForceText=$false},
    @{Cell="D31"; Value='0.05626'; ForceText=$true},
    @{Cell="E31"; Value='  -1.38%  '; ForceText=$false},
    @{Cell="D32"; Value='4.116'; ForceText=$true},
    @{Cell="E32"; Value='  -0.36%  '; ForceText=$false},
    @{Cell="D33"; Value='4.046'; ForceText=$true},
    @{Cell="E33"; Value='  +0.64%  '; ForceText=$false},
    @{Cell="D34"; Value='1.851'; ForceText=$true},
    @{Cell="E34"; Value='  +0.26%  '; ForceText=$false},
    @{Cell="E35"; Value='  +0.11%  '; ForceText=$false},
    @{Cell="D36"; Value='0.7095'; ForceText=$true},
    @{Cell="E36"; Value='  -0.92%  '; ForceText=$false},
    @{Cell="D37"; Value='2.585'; ForceText=$true},
    @{Cell="E37"; Value='  -0.19%  '; ForceText=$false},
    @{Cell="D38"; Value='1.234.70'; ForceText=$false},
    @{Cell="E38"; Value='  -1.38%  '; ForceText=$false},
    @{Cell="D39"; Value='2.775'; ForceText=$true},
    @{Cell="E39"; Value='  -0.24%  '; ForceText=$false},
    @{Cell="E40"; Value='  -1.04%  '; ForceText=$false},
    @{Cell="D41"; Value='6.477'; ForceText=$true},
    @{Cell="E41"; Value='  +4.69%  '; ForceText=$false},
    @{Cell="D42"; Value='0.9072'; ForceText=$true},
    @{Cell="E42"; Value='  -0.61%  '; ForceText=$false},
    @{Cell="D43"; Value='1.000'; ForceText=$true},
    @{Cell="E43"; Value='  +0.00%  '; ForceText=$false},
    @{Cell="D44"; Value='1.995.58'; ForceText=$false},
    @{Cell="E44"; Value='  -0.32%  '; ForceText=$false},
    @{Cell="D45"; Value='101.38'; ForceText=$true},
    @{Cell="E45"; Value='  -0.46%  '; ForceText=$false},
    @{Cell="D46"; Value='66.07'; ForceText=$true},
    @{Cell="E46"; Value='  -0.10%  '; ForceText=$false},
    @{Cell="E47"; Value='  +3.06%  '; ForceText=$false},
    @{Cell="D48"; Value='7.166'; ForceText=$true},
    @{Cell="E48"; Value='  +1.34%  '; ForceText=$false},
    @{Cell="D49"; Value='0.4013'; ForceText=$true},
    @{Cell="E49"; Value='  -0.40%  '; ForceText=$false},
    @{Cell="B50"; Value='RenderToken'; ForceText=$false},
    @{Cell="C50"; Value='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; ForceText=$false},
    @{Cell="D50"; Value='1.687'; ForceText=$true},
    @{Cell="E50"; Value='  -0.83%  '; ForceText=$false},
    @{Cell="B51"; Value='EnergySwap'; ForceText=$false},
    @{Cell="C51"; Value='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; ForceText=$false},
    @{Cell="D51"; Value='8.993'; ForceText=$true},
    @{Cell="E51"; Value='  -1.36%  '; ForceText=$false}
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Force text storage so numeric-looking strings (e.g. "1.000") keep their
        # exact literal formatting instead of being auto-coerced to a Double by Excel.
        $origStyle = $r.Style
        $r.NumberFormat = "@"
        $r.Value = $u.Value
        $r.Style = $origStyle
    } else {
        $r.Value = $u.Value
    }
}
